# Update countries & provincias Spain
#
# The "Pais" sheet lists countries sorted by Column B (Casos totales)
# descending. This update refreshes several countries' daily figures;
# because a few countries' totals crossed a neighbour's total, those rows
# swap places (the shared-string table gets reshuffled as a side effect,
# but all we need to do is write the correct country name + figures into
# each affected row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("A$row").Value = $country
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# row, country,              Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
Set-Row 4   "Estados Unidos"      675640 27637 57271 583847 13369 2079 34522
Set-Row 8   "Alemania"            137698  2945 77000  56646  4288  248  4052

Set-Row 21  "India"                13430  1060  1768  11214     0   26   448
Set-Row 22  "Irlanda"              13271   724    77  12708   156   42   486

Set-Row 33  "Noruega"               6896    99    32   6712    64    2   152
Set-Row 34  "Dinamarca"             6879   198  3023   3535    92   12   321

Set-Row 92  "Costa Rica"             642    16    74    564    11    0     4

Set-Row 96  "Burkina Faso"           546     4   257    257     0    0    32
Set-Row 98  "Uruguay"                502     9   286    207    13    0     9

Set-Row 113 "Montenegro"             303    15    55    244     7    0     4

Set-Row 172 "Maldivas"                25     3    16      9     0    0     0
Set-Row 173 "Zimbabue"                23     0     1     19     0    0     3
Set-Row 174 "Antigua y Barbuda"       23     0     3     17     1    1     3
